# Refresh the Price (D) / Volume(1h) (E) columns with the latest scrape values.
# Matches the GitHub Actions "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several Price values look numeric (e.g. "218.10", "0.07740"); the sheet
# stores them as plain text, so force Text format before assigning, then
# drop back to the Normal style so no stray number-format index remains.
$textPriceCells = @("D4","D5","D6","D10","D11","D12","D14","D16","D18","D19","D20","D22","D24","D25","D26","D27","D29","D30","D31","D33","D35","D36","D37","D38","D39","D40","D42","D45","D47","D48","D49","D50","D51")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.090.33"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.650.72"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "218.10"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "0.5281"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "20.35"
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "4.467"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "1.650.44"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "0.5448"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "0.0₅8109"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "65.14"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "26.100.73"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "4.540"
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("D20").Value = "193.63"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "5.973"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "140.02"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").Value = "0.1239"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "7.238"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").Value = "16.14"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "0.05901"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "1.280"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "3.500"
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("D33").Value = "1.547"
$ws.Range("E33").Value = "  -5.91%  "
$ws.Range("D35").Value = "0.9420"
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("D36").Value = "2.760"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "0.5654"
$ws.Range("E37").Value = "  -4.18%  "
$ws.Range("D38").Value = "0.01601"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "5.841"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "0.8442"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "100.74"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "1.006.64"
$ws.Range("E43").Value = "  -3.07%  "
$ws.Range("D44").Value = "1.798.07"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "56.79"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").Value = "0.4288"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").Value = "1.475"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").Value = "0.05150"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").Value = "7.807"
$ws.Range("E51").Value = "  -3.72%  "

foreach ($addr in $textPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
